$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement the "剩余" (remaining) value in column E by 1 for every data row
# (rows 2 through 99), except row 36 whose start-date value looks malformed
# and was left untouched in the source edit.
for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value()
    $cell.Value = $current - 1
}
